# "Running all the suites" — flip the Runmode column (C) to "Y" for every
# test case row on the "Test Cases" sheet, and update the active selection
# to reflect that column (C2:C18, active cell C2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Column C ("Runmode") -> "Y" for every data row (row 2 through row 18).
$ws.Range("C2:C18").Value = "Y"

# Match the selection recorded in the saved workbook.
$ws.Range("C2:C18").Select()
